$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Ключ"
$ws.Range("B1").Value = "Значение на русском"
$ws.Range("C1").Value = "Значение на казахском"

$ws.Range("B11").Select()
